$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "ROGER"
$ws.Range("C4").Value = 200
$ws.Range("D4").Value = 376
$ws.Range("E4").Value = 376

# Row 5
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "WAGNER"
$ws.Range("C5").Value = 300
$ws.Range("D5").Value = 556
$ws.Range("E5").Value = 556

# Row 6
$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "RODRIGO"
$ws.Range("C6").Value = 70
$ws.Range("D6").Value = 298
$ws.Range("E6").Value = 298

# Row 7
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "THIAGO FERNANDES RODRIGUES"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 296
$ws.Range("E7").Value = 1296
$ws.Range("F7").Value = "1234ab"
